{"js": "// Replace the 100 math-expression cell values in the single 20x5 table,\n// in row-major order, matching the source diff (old -> new per cell).\nconst newValues = [\n  [\"26+64=\", \"3+10=\", \"70-23=\", \"16+69=\", \"48+50=\"],\n  [\"65+4=\", \"66+18=\", \"88-35=\", \"48+35=\", \"84-81=\"],\n  [\"85-75=\", \"54+4=\", \"35+9=\", \"69-9=\", \"58-5=\"],\n  [\"27+59=\", \"45-28=\", \"16+4=\", \"4+42=\", \"66+23=\"],\n  [\"13-5=\", \"24+12=\", \"35-33=\", \"94-46=\", \"53+25=\"],\n  [\"66-7=\", \"10+10=\", \"44+41=\", \"62-51=\", \"42-22=\"],\n  [\"46-25=\", \"31+26=\", \"40-0=\", \"23-6=\", \"50-3=\"],\n  [\"42-14=\", \"89-8=\", \"42+45=\", \"22+64=\", \"34-13=\"],\n  [\"94-32=\", \"18+80=\", \"16+41=\", \"96-37=\", \"30-27=\"],\n  [\"27+2=\", \"97-45=\", \"36+35=\", \"95-3=\", \"4+75=\"],\n  [\"92-79=\", \"24+3=\", \"90-87=\", \"35-4=\", \"89+7=\"],\n  [\"49+39=\", \"86+13=\", \"54+5=\", \"57+27=\", \"74+0=\"],\n  [\"89-3=\", \"83-83=\", \"91-49=\", \"15+78=\", \"24-12=\"],\n  [\"19+34=\", \"68-10=\", \"70-16=\", \"78-60=\", \"72-61=\"],\n  [\"70-16=\", \"62-49=\", \"99-20=\", \"65-7=\", \"65+19=\"],\n  [\"67-54=\", \"31+37=\", \"66+9=\", \"18+76=\", \"57+0=\"],\n  [\"31+29=\", \"26+11=\", \"35+18=\", \"70-1=\", \"4+58=\"],\n  [\"21+11=\", \"16+43=\", \"22+8=\", \"63+11=\", \"82-42=\"],\n  [\"71-23=\", \"44-29=\", \"49+32=\", \"72-69=\", \"77-2=\"],\n  [\"84-69=\", \"13+40=\", \"75+14=\", \"26+25=\", \"2+35=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst existingCols = table.values.length > 0 ? table.values[0].length : 0;\nif (table.rowCount !== newValues.length || existingCols !== newValues[0].length) {\n  throw new Error(\n    `Unexpected table shape: ${table.rowCount}x${existingCols}, expected ${newValues.length}x${newValues[0].length}`\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the 100 math-expression cells in the single table (20 rows x 5 cols)\n# with new values, in row-major order, matching the source diff.\n$newValues = @(\n    \"26+64=\",\n    \"3+10=\",\n    \"70-23=\",\n    \"16+69=\",\n    \"48+50=\",\n    \"65+4=\",\n    \"66+18=\",\n    \"88-35=\",\n    \"48+35=\",\n    \"84-81=\",\n    \"85-75=\",\n    \"54+4=\",\n    \"35+9=\",\n    \"69-9=\",\n    \"58-5=\",\n    \"27+59=\",\n    \"45-28=\",\n    \"16+4=\",\n    \"4+42=\",\n    \"66+23=\",\n    \"13-5=\",\n    \"24+12=\",\n    \"35-33=\",\n    \"94-46=\",\n    \"53+25=\",\n    \"66-7=\",\n    \"10+10=\",\n    \"44+41=\",\n    \"62-51=\",\n    \"42-22=\",\n    \"46-25=\",\n    \"31+26=\",\n    \"40-0=\",\n    \"23-6=\",\n    \"50-3=\",\n    \"42-14=\",\n    \"89-8=\",\n    \"42+45=\",\n    \"22+64=\",\n    \"34-13=\",\n    \"94-32=\",\n    \"18+80=\",\n    \"16+41=\",\n    \"96-37=\",\n    \"30-27=\",\n    \"27+2=\",\n    \"97-45=\",\n    \"36+35=\",\n    \"95-3=\",\n    \"4+75=\",\n    \"92-79=\",\n    \"24+3=\",\n    \"90-87=\",\n    \"35-4=\",\n    \"89+7=\",\n    \"49+39=\",\n    \"86+13=\",\n    \"54+5=\",\n    \"57+27=\",\n    \"74+0=\",\n    \"89-3=\",\n    \"83-83=\",\n    \"91-49=\",\n    \"15+78=\",\n    \"24-12=\",\n    \"19+34=\",\n    \"68-10=\",\n    \"70-16=\",\n    \"78-60=\",\n    \"72-61=\",\n    \"70-16=\",\n    \"62-49=\",\n    \"99-20=\",\n    \"65-7=\",\n    \"65+19=\",\n    \"67-54=\",\n    \"31+37=\",\n    \"66+9=\",\n    \"18+76=\",\n    \"57+0=\",\n    \"31+29=\",\n    \"26+11=\",\n    \"35+18=\",\n    \"70-1=\",\n    \"4+58=\",\n    \"21+11=\",\n    \"16+43=\",\n    \"22+8=\",\n    \"63+11=\",\n    \"82-42=\",\n    \"71-23=\",\n    \"44-29=\",\n    \"49+32=\",\n    \"72-69=\",\n    \"77-2=\",\n    \"84-69=\",\n    \"13+40=\",\n    \"75+14=\",\n    \"26+25=\",\n    \"2+35=\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\nif (($rows * $cols) -ne $newValues.Count) {\n    throw \"Unexpected table shape: $rows x $cols (= $($rows*$cols) cells), expected $($newValues.Count) cells\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n\nWrite-Output (\"Updated \" + $i + \" cells\")\n"}
